# Commit: "Add article type, game and mod controler. New design. New jquery,
# Signed-off-by:Igor Peshkov <igor.peshkov@gmail.com>"
#
# Turns the two previously-blank backlog rows (9 & 10) on the "Tasks" sheet
# into real task rows, and rewrites row 8's task from the old
# "restore password" entry into the new "article type controller" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- Row 8: replace the old "restore password" task with the new
#     "article type" task (Task / Task Location columns only; Priority
#     "Normal" and Status "Opened" stay as they already were). The long
#     Task description (B8) is filled in last, below, after rows 9-10. ---
$ws.Range("C8").Value = "ArticleTypeController"

# --- Row 9: was completely empty -> new "cancel button" task. ---
$ws.Range("B9").Value = "Добавить кнопки `"отмена`" на страницы редактирования статей."
$ws.Range("C9").Value = "ArctileController/Edit"
$ws.Range("D9").Value = "Low"
$ws.Range("F9").Value = "После нажатия на эту кнопку следует возратить на страницу просмотра статьи"

# --- Row 10: was completely empty -> new "ACL class" task. ---
$ws.Range("B10").Value = "Переделать ACL класс. Заменить два запроса достающие сначала номер роли, а потом имя роли на один с джойном."
$ws.Range("C10").Value = "ACL class"
$ws.Range("D10").Value = "Mdium"

# Rows 9 & 10 now have a Status too ("Opened"), matching the red "Bad"
# badge style already used by row 8's Status cell (E8) -- copy that exact
# formatting across instead of building a brand-new style entry.
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E9").Value = "Opened"
$ws.Range("E10").Value = "Opened"

# Row 8's long Task description, filled in last.
$ws.Range("B8").Value = "Реализовать добавления типов статей через админку в бд. А так же их выбор при создании статьи."

# The longer wrapped text in B/C now needs two lines, so these three rows
# grow from the default 15pt to 30pt.
$ws.Range("A8:A10").EntireRow.RowHeight = 30

# Matches the saved selection recorded in the sheet view.
$ws.Range("B9").Select()
